# Refresh the hourly cryptos snapshot (price + 1h volume change), and
# account for the new "BabyDogeCoin" entry that pushed RocketPoolETH and
# Mantle down a row.
#
# Price values in column D are free-form text (e.g. "28.934.88",
# "0.000008313") rather than real numbers, so each one is written with a
# leading single-quote to force Excel to keep it as text instead of
# reinterpreting it (and mangling precision/format) as a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Price($row, $price) {
    $ws.Range("D$row").Value = "'" + $price
}

function Set-Volume($row, $volume) {
    $ws.Range("E$row").Value = $volume
}

# Row 2 - Bitcoin
Set-Price 2 "28.934.88"
Set-Volume 2 "  -1.51%  "

# Row 3 - Ethereum
Set-Price 3 "1.832.72"
Set-Volume 3 "  -1.88%  "

# Row 4 - TetherUSD
Set-Price 4 "1.0000"
Set-Volume 4 "  -0.12%  "

# Row 5 - BNB
Set-Price 5 "245.20"
Set-Volume 5 "  +0.60%  "

# Row 6 - XRP
Set-Price 6 "0.6904"
Set-Volume 6 "  -1.87%  "

# Row 7 - USDC
Set-Volume 7 "  -0.08%  "

# Row 8 - Dogecoin
Set-Price 8 "0.07697"
Set-Volume 8 "  -3.03%  "

# Row 9 - Cardano
Set-Volume 9 "  -2.60%  "

# Row 10 - Solana
Set-Price 10 "23.48"
Set-Volume 10 "  -3.97%  "

# Row 11 - TRON
Set-Price 11 "0.07813"
Set-Volume 11 "  -0.35%  "

# Row 12 - WrappedEther
Set-Price 12 "1.836.43"
Set-Volume 12 "  -1.60%  "

# Row 13 - Polkadot
Set-Price 13 "5.084"
Set-Volume 13 "  -1.86%  "

# Row 14 - Litecoin
Set-Volume 14 "  -3.70%  "

# Row 15 - Polygon
Set-Price 15 "0.6803"
Set-Volume 15 "  -2.89%  "

# Row 16 - Uniswap
Set-Price 16 "6.422"
Set-Volume 16 "  -1.68%  "

# Row 17 - ShibaInu
Set-Price 17 "0.000008313"
Set-Volume 17 "  -0.94%  "

# Row 18 - WrappedBTC
Set-Price 18 "28.921.16"
Set-Volume 18 "  -1.53%  "

# Row 19 - BitcoinCash
Set-Price 19 "243.25"
Set-Volume 19 "  -4.03%  "

# Row 20 - WrappedliquidstakedEther2.0
Set-Price 20 "2.081.56"
Set-Volume 20 "  -1.34%  "

# Row 21 - Avalanche
Set-Volume 21 "  -2.87%  "

# Row 22 - Dai
Set-Volume 22 "  -0.01%  "

# Row 23 - Chainlink
Set-Price 23 "7.469"
Set-Volume 23 "  -2.42%  "

# Row 25 - Monero
Set-Price 25 "163.07"
Set-Volume 25 "  +0.82%  "

# Row 26 - Stellar
Set-Price 26 "0.1467"
Set-Volume 26 "  -5.51%  "

# Row 27 - Cosmos
Set-Price 27 "8.798"
Set-Volume 27 "  -2.31%  "

# Row 28 - EthereumClassic
Set-Volume 28 "  -3.42%  "

# Row 29 - PancakeSwap
Set-Volume 29 "  +3.18%  "

# Row 30 - Filecoin
Set-Price 30 "4.213"

# Row 31 - InternetComputer(DFINITY)
Set-Price 31 "4.152"
Set-Volume 31 "  -2.28%  "

# Row 32 - Toncoin
Set-Price 32 "1.176"
Set-Volume 32 "  -3.15%  "

# Row 33 - Hedera
Set-Price 33 "0.05113"
Set-Volume 33 "  -3.14%  "

# Row 34 - ImmutableX
Set-Price 34 "0.7684"
Set-Volume 34 "  +2.77%  "

# Row 35 - LidoDAOToken
Set-Price 35 "1.839"
Set-Volume 35 "  -3.00%  "

# Row 36 - ARBITRUM
Set-Price 36 "1.144"
Set-Volume 36 "  -2.64%  "

# Row 37 - HuobiToken
Set-Price 37 "2.686"

# Row 38 - VeChain
Set-Price 38 "0.01850"
Set-Volume 38 "  -1.57%  "

# Row 39 - Maker
Set-Price 39 "1.238.47"
Set-Volume 39 "  -2.94%  "

# Row 40 - MXToken
Set-Price 40 "2.697"

# Row 41 - TrustWalletToken
Set-Price 41 "0.9366"
Set-Volume 41 "  +5.06%  "

# Row 42 - Quant
Set-Price 42 "107.79"
Set-Volume 42 "  -0.68%  "

# Row 43 - PaxDollar
Set-Price 43 "0.9992"
Set-Volume 43 "  -0.17%  "

# Row 44 - FraxShare
Set-Volume 44 "  -5.39%  "

# Row 45 - EnergySwap
Set-Price 45 "9.574"
Set-Volume 45 "  -0.25%  "

# Row 46 - now BabyDogeCoin (new entrant, bumping the two rows below down)
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-Price 46 "0.00000000122"
Set-Volume 46 "  -4.13%  "

# Row 47 - now RocketPoolETH (was row 46)
$ws.Range("B47").Value = "RocketPoolETH"
$ws.Range("C47").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
Set-Price 47 "1.981.43"
Set-Volume 47 "  -1.62%  "

# Row 48 - now Mantle (was row 47)
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-Price 48 "0.5173"
Set-Volume 48 "  -0.14%  "

# Row 49 - Aave
Set-Price 49 "64.55"
Set-Volume 49 "  -9.19%  "

# Row 50 - RenderToken
Set-Price 50 "1.750"
Set-Volume 50 "  -2.69%  "

# Row 51 - TheSandbox
Set-Price 51 "0.4195"
Set-Volume 51 "  -2.46%  "
